$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = "iteration 1"

$ws.Range("B3").Value = 461286
$ws.Range("C3").Value = 136725
$ws.Range("D3").Value = 173623
$ws.Range("E3").Value = 75572
$ws.Range("F3").Value = 75366
